# StructureDefinition-enrollment-type.xlsx update
# - Version bump 5.0.0 -> 6.0.0
# - Date bump
# - Publisher filled in with "Alvearie Team"
# - Second "Contact" row replaced with "Jurisdiction" / "United States of America"
# - Duplicate "Contact" row removed
# - Elements sheet root row's Short/Definition updated to resource-specific text

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Enrollment Type"
$elements.Range("L2").Value = "Code for the type of enrollment"
